$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.996.77'
$ws.Range("E2").Value = '  +0.26%  '

$ws.Range("D3").Value = '2.368.13'
$ws.Range("E3").Value = '  -0.70%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").Value = '315.73'
$ws.Range("E5").Value = '  -3.42%  '

$ws.Range("D6").Value = '108.94'
$ws.Range("E6").Value = '  +9.43%  '

$ws.Range("D7").Value = '0.639'
$ws.Range("E7").Value = '  +0.31%  '

$ws.Range("E8").Value = '  -0.14%  '

$ws.Range("D9").Value = '0.623'
$ws.Range("E9").Value = '  -0.43%  '

$ws.Range("D10").Value = '41.38'
$ws.Range("E10").Value = '  +3.57%  '

$ws.Range("D11").Value = '0.0931'
$ws.Range("E11").Value = '  +0.92%  '

$ws.Range("D12").Value = '8.64'
$ws.Range("E12").Value = '  +2.95%  '

$ws.Range("E13").Value = '  -0.21%  '

$ws.Range("E14").Value = '  +1.35%  '

$ws.Range("D15").Value = '15.98'
$ws.Range("E15").Value = '  -3.15%  '

$ws.Range("D16").Value = '2.724.24'
$ws.Range("E16").Value = '  -0.85%  '

$ws.Range("D17").Value = '2.349.89'
$ws.Range("E17").Value = '  -1.61%  '

$ws.Range("D18").Value = '42.941.04'
$ws.Range("E18").Value = '  +0.20%  '

$ws.Range("D19").Value = '7.68'
$ws.Range("E19").Value = '  -1.43%  '

$ws.Range("E20").Value = '  +0.57%  '

$ws.Range("D21").Value = '76.43'
$ws.Range("E21").Value = '  +1.42%  '

$ws.Range("D22").Value = '3.59'
$ws.Range("E22").Value = '  -4.43%  '

$ws.Range("D23").Value = '268.82'
$ws.Range("E23").Value = '  -1.55%  '

$ws.Range("D24").Value = '2.35'
$ws.Range("E24").Value = '  +0.41%  '

$ws.Range("D25").Value = '9.54'
$ws.Range("E25").Value = '  -6.23%  '

$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("D27").Value = '11.51'
$ws.Range("E27").Value = '  +0.22%  '

$ws.Range("D28").Value = '23.45'
$ws.Range("E28").Value = '  -1.78%  '

$ws.Range("E29").Value = '  +2.03%  '

$ws.Range("D30").Value = '36.93'
$ws.Range("E30").Value = '  +4.70%  '

$ws.Range("D31").Value = '168.89'
$ws.Range("E31").Value = '  -2.27%  '

$ws.Range("E32").Value = '  +0.16%  '

$ws.Range("E33").Value = '  +2.58%  '

$ws.Range("E34").Value = '  -6.58%  '

$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").Value = '0.121'
$ws.Range("E35").Value = '  +14.61%  '

$ws.Range("B36").Value = 'Stellar'
$ws.Range("C36").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D36").Value = '0.132'
$ws.Range("E36").Value = '  -0.40%  '

$ws.Range("E37").Value = '  +2.17%  '

$ws.Range("E38").Value = '  +1.43%  '

$ws.Range("D39").Value = '3.88'
$ws.Range("E39").Value = '  +0.67%  '

$ws.Range("D40").Value = '2.69'
$ws.Range("E40").Value = '  -6.01%  '

$ws.Range("D41").Value = '105.23'
$ws.Range("E41").Value = '  +10.26%  '

$ws.Range("D42").Value = '1.52'
$ws.Range("E42").Value = '  -0.72%  '

$ws.Range("D43").Value = '0.239'
$ws.Range("E43").Value = '  +4.80%  '

$ws.Range("D44").Value = '71.54'
$ws.Range("E44").Value = '  +3.90%  '

$ws.Range("E45").Value = '  -0.07%  '

$ws.Range("D46").Value = '12.66'
$ws.Range("E46").Value = '  +6.62%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '114.57'
$ws.Range("E47").Value = '  -1.68%  '

$ws.Range("B48").Value = 'ordi'
$ws.Range("C48").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D48").Value = '81.02'
$ws.Range("E48").Value = '  +19.68%  '

$ws.Range("D49").Value = '5.59'
$ws.Range("E49").Value = '  +2.78%  '

$ws.Range("D50").Value = '9.19'
$ws.Range("E50").Value = '  +2.00%  '

$ws.Range("E51").Value = '  +2.44%  '
